$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(177,106,125,114,143,236,242,134,115),
    @(0,0,2,4,5,6,3,4,2),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,4,4,4,0,3,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,1,1,1,0,1),
    @(0,0,2,0,0,1,2,1,1),
    @(0,0,0,0,0,0,0,0,0),
    @(51,26,42,39,24,56,53,31,44),
    @(9,8,7,9,0,9,8,9,8),
    @(13,2,11,9,7,9,6,4,4),
    @(3,3,4,1,2,5,6,0,4),
    @(9,4,8,8,6,13,10,7,14),
    @(7,6,8,0,6,13,16,6,14),
    @(10,3,4,12,3,7,7,5,0),
    @(25,12,24,6,7,19,21,26,13),
    @(0,0,0,0,0,0,0,0,0),
    @(11,3,12,1,3,15,9,14,5),
    @(0,0,0,0,0,0,0,0,0),
    @(7,8,7,2,1,3,5,2,4),
    @(7,1,5,3,3,1,7,10,4),
    @(36,16,15,15,29,42,40,13,13),
    @(9,4,5,2,13,6,12,1,4),
    @(4,7,0,7,7,9,3,5,2),
    @(10,3,3,6,8,14,8,2,3),
    @(13,2,7,0,1,13,17,5,4),
    @(9,15,16,14,34,27,35,21,24),
    @(3,9,13,4,12,13,17,16,16),
    @(2,5,3,2,7,5,3,1,5),
    @(0,0,0,0,0,0,0,0,0),
    @(4,1,0,8,15,9,15,4,3),
    @(0,0,0,0,0,0,0,0,0),
    @(31,28,18,30,36,49,50,21,8),
    @(9,9,9,9,9,9,9,7,3),
    @(0,0,0,0,0,0,0,0,0),
    @(12,10,6,15,18,22,23,6,3),
    @(4,3,3,3,4,6,5,4,2),
    @(6,6,0,3,5,12,13,4,0),
    @(25,9,8,6,8,37,40,18,11),
    @(10,7,6,5,6,17,13,11,8),
    @(7,1,2,1,0,10,8,7,2),
    @(5,0,0,0,0,1,8,0,0),
    @(3,1,0,0,2,9,11,0,1),
    @(0,0,0,0,0,0,0,0,0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $rowValues[$c]
    }
}
